$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new "Save" column (H) to the s_vals sheet, mirroring the existing
# header formatting (bold/centered/bordered) used by the other headers
# (TB, d2S, K, IP, Win, sum) in row 1, and a numeric 0 in row 2.

# Copy G1's formatting (and value, temporarily) into H1 so the new header
# cell picks up the same style as the rest of row 1, then overwrite the
# value with the real header text.
$ws.Range("G1").Copy($ws.Range("H1"))
$ws.Range("H1").Value = "Save"

# Data row: plain numeric 0, same as the other numeric cells in row 2.
$ws.Range("H2").Value = 0
